# Added logger information and comments
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("addInput")

# Remove the extra "abhai" logger/output column (column D) that is no
# longer needed now that logger information is captured elsewhere.
$ws.Columns.Item(4).Delete() | Out-Null

# Update the sample add-input row with the new first/second numbers and
# their resulting output, used as comment/demo data for the logger.
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 12

# Move the active selection to reflect where the user last left off.
$ws.Range("C11").Select() | Out-Null
